$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.226.46"
$ws.Range("E2").Value = "  -0.63%  "
$ws.Range("D3").Value = "2.267.40"
$ws.Range("E3").Value = "  -1.01%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "300.19"
$ws.Range("E5").Value = "  -0.66%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "95.84"
$ws.Range("E6").Value = "  -2.39%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.494"
$ws.Range("E7").Value = "  -2.20%  "
$ws.Range("E8").Value = "  +0.07%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.491"
$ws.Range("E9").Value = "  -1.56%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "33.11"
$ws.Range("E10").Value = "  -3.56%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0787"
$ws.Range("E11").Value = "  -0.02%  "
$ws.Range("E12").Value = "  -6.38%  "
$ws.Range("E13").Value = "  +0.84%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.65"
$ws.Range("E14").Value = "  -0.87%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.67"
$ws.Range("E15").Value = "  +0.40%  "
$ws.Range("D16").Value = "2.619.85"
$ws.Range("E16").Value = "  -1.00%  "
$ws.Range("D17").Value = "2.270.99"
$ws.Range("E17").Value = "  +0.04%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.782"
$ws.Range("E18").Value = "  -2.45%  "
$ws.Range("D19").Value = "42.150.72"
$ws.Range("E19").Value = "  -0.63%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.63"
$ws.Range("E20").Value = "  +1.31%  "
$ws.Range("D21").Value = "0.0₃0888"
$ws.Range("E21").Value = "  -1.11%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.98"
$ws.Range("E22").Value = "  -0.66%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "66.25"
$ws.Range("E23").Value = "  -2.13%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "234.86"
$ws.Range("E24").Value = "  +0.00%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.96"
$ws.Range("E25").Value = "  -0.38%  "
$ws.Range("B26").Value = "Dai"
$ws.Range("C26").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("E26").Value = "  +0.15%  "
$ws.Range("B27").Value = "PancakeSwap"
$ws.Range("C27").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.46"
$ws.Range("E27").Value = "  -2.40%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "23.88"
$ws.Range("E28").Value = "  -4.36%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.29"
$ws.Range("E29").Value = "  +0.16%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "168.03"
$ws.Range("E30").Value = "  +3.05%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "9.15"
$ws.Range("E31").Value = "  +0.01%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "33.60"
$ws.Range("E32").Value = "  -3.00%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.00"
$ws.Range("E33").Value = "  -0.01%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.87"
$ws.Range("E34").Value = "  -2.52%  "
$ws.Range("E35").Value = "  -0.43%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "16.59"
$ws.Range("E36").Value = "  -1.99%  "
$ws.Range("E37").Value = "  -3.32%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0684"
$ws.Range("E38").Value = "  -3.84%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.77"
$ws.Range("E39").Value = "  -3.71%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0984"
$ws.Range("E40").Value = "  -2.11%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.108"
$ws.Range("E41").Value = "  -2.46%  "
$ws.Range("E42").Value = "  -4.17%  "
$ws.Range("E43").Value = "  -0.19%  "
$ws.Range("D44").Value = "1.970.29"
$ws.Range("E44").Value = "  -0.12%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0277"
$ws.Range("E45").Value = "  -0.86%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "17.42"
$ws.Range("E46").Value = "  -6.10%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.53"
$ws.Range("E47").Value = "  -5.95%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.77"
$ws.Range("E48").Value = "  -4.46%  "
$ws.Range("D49").Value = "2.492.95"
$ws.Range("E49").Value = "  -0.98%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "52.39"
$ws.Range("E50").Value = "  -5.41%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.47"
$ws.Range("E51").Value = "  -0.35%  "
